$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '44.169.56'
$ws.Range('E2').Value = '  +5.38%  '
$ws.Range('D3').Value = '2.262.66'
$ws.Range('E3').Value = '  +2.50%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '230.52'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.31%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.632'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +2.56%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '63.58'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +4.87%  '
$ws.Range('E8').Value = '  +0.08%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.448'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +11.63%  '
$ws.Range('E10').Value = '  +14.62%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '56.83'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -0.68%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '26.36'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +19.89%  '
$ws.Range('E13').Value = '  +2.47%  '
$ws.Range('D14').Value = '2.597.96'
$ws.Range('E14').Value = '  +2.37%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '15.70'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +1.97%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '6.06'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +8.84%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.834'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +4.67%  '
$ws.Range('D18').Value = '2.264.19'
$ws.Range('E18').Value = '  +2.12%  '
$ws.Range('D19').Value = '43.947.92'
$ws.Range('E19').Value = '  +4.94%  '
$ws.Range('E20').Value = '  +7.45%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '73.60'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +2.06%  '
$ws.Range('E22').Value = '  -0.72%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '255.08'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +5.18%  '
$ws.Range('E24').Value = '  +0.05%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.42'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +3.72%  '
$ws.Range('E26').Value = '  -7.09%  '
$ws.Range('B27').Value = 'WEMIXToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '3.36'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +26.81%  '
$ws.Range('B28').Value = 'Cosmos'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '10.12'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +5.25%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '171.98'
$ws.Range('D29').Style = "Normal"
$ws.Range('B30').Value = 'Kaspa'
$ws.Range('C30').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.137'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -1.62%  '
$ws.Range('B31').Value = 'EthereumClassic'
$ws.Range('C31').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '20.76'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +2.68%  '
$ws.Range('E32').Value = '  -2.49%  '
$ws.Range('E33').Value = '  +3.04%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.0678'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +5.00%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '4.76'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +3.68%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '4.85'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -1.81%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '3.82'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +8.19%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '6.71'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +7.00%  '
$ws.Range('E39').Value = '  -0.18%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.0257'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +5.73%  '
$ws.Range('E41').Value = '  -0.11%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '17.64'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +9.92%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '8.30'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -2.64%  '
$ws.Range('E44').Value = '  +1.46%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '97.63'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +1.14%  '
$ws.Range('B46').Value = 'TrustWalletToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '1.19'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +0.02%  '
$ws.Range('B47').Value = 'FTXToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '4.39'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +0.75%  '
$ws.Range('E48').Value = '  -4.45%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '10.07'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +19.58%  '
$ws.Range('D50').Value = '1.446.00'
$ws.Range('E50').Value = '  -0.57%  '
$ws.Range('E51').Value = '  +4.36%  '
